$wb = $excel.ActiveWorkbook

# Sheet "展览" (sheet1) - first worksheet
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("F2").Value = 6506
$ws1.Range("F7").Value = 8
$ws1.Range("F8").Value = 536
$ws1.Range("F12").Value = 161
$ws1.Range("F13").Value = 382
$ws1.Range("F14").Value = 953
$ws1.Range("F15").Value = 3209
$ws1.Range("F18").Value = 1874

# Sheet "全部类型" (sheet4) - fourth worksheet
$ws4 = $wb.Worksheets.Item(4)
$ws4.Range("F2").Value = 6506
$ws4.Range("F7").Value = 8
$ws4.Range("F9").Value = 536
$ws4.Range("F13").Value = 161
$ws4.Range("F14").Value = 382
$ws4.Range("F15").Value = 953
$ws4.Range("F16").Value = 3209
$ws4.Range("F19").Value = 1874
